# Applies the "regression_hdi/hdi_2015.xlsx" revision:
#  - renames the two worksheets (Sheet2 -> raw_data, Sheet1 -> variables)
#  - removes the stray row 30 (a lone space character) on the variables sheet
#  - adds a new "ignore_variables" column to Table2 on the variables sheet,
#    flagging Country / HDI Rank / HDI / HDI rank2 as variables to ignore
#  - tidies up the selections left behind by the edit

$wb = $excel.ActiveWorkbook

$rawData = $wb.Worksheets.Item("Sheet2")
$variables = $wb.Worksheets.Item("Sheet1")

# --- rename sheets ---------------------------------------------------
$rawData.Name = "raw_data"
$variables.Name = "variables"

# --- drop the leftover row with a single blank-space string ----------
$variables.Rows("30:30").Delete()

# --- extend Table2 with the new "ignore_variables" column -------------
$table = $variables.ListObjects.Item("Table2")
$newColumn = $table.ListColumns.Add()
$newColumn.Range.Cells(1, 1).Value = "ignore_variables"

# match the header cell's look to the rest of the header row
$variables.Range("B1").Copy()
$variables.Range("E1").PasteSpecial(-4122)

# match the new data cells' look to the rest of the table body, but
# without the thin left/right borders the inner columns have
$variables.Range("B2").Copy()
$dataRange = $variables.Range("E2:E25")
$dataRange.PasteSpecial(-4122)
$dataRange.Borders.LineStyle = -4142

# Flag Country, HDI Rank, HDI (Human Development Index) and HDI rank2
# as variables that should be ignored; everything else stays blank.
$variables.Range("E2").Value = 1
$variables.Range("E3").Value = 1
$variables.Range("E4").Value = 1
$variables.Range("E17").Value = 1

# --- restore selections -----------------------------------------------
$rawData.Range("E23").Select() | Out-Null
$variables.Range("E1").Select() | Out-Null
